# Weekly update: a new daily record for "Perejil" (parsley) at Terminal La
# Palmera de La Serena is inserted as the new row 18. All the existing
# historical rows (old rows 18-79) shift down by one row (to rows 19-80);
# this also extends the used range from A1:R79 to A1:R80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 18, pushing every row below it down.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(18, 1).Value  = 8
$ws.Cells.Item(18, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(18, 3).Value  = "Coquimbo"
$ws.Cells.Item(18, 4).Value  = 44453
$ws.Cells.Item(18, 5).Value  = 4
$ws.Cells.Item(18, 6).Value  = 100112044
$ws.Cells.Item(18, 7).Value  = "Perejil"
$ws.Cells.Item(18, 8).Value  = "Sin especificar"
$ws.Cells.Item(18, 9).Value  = "Primera"
$ws.Cells.Item(18, 10).Value = 3200
$ws.Cells.Item(18, 11).Value = 2000
$ws.Cells.Item(18, 12).Value = 2500
$ws.Cells.Item(18, 13).Value = 2250
$ws.Cells.Item(18, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(18, 16).Value = 1500
$ws.Cells.Item(18, 17).Value = 1.5
$ws.Cells.Item(18, 18).Value = "Hortaliza"
